## Add a new "RTM" (Requirement Traceability Matrix) worksheet at the end
## of the workbook, populate it with data, format it as a table, and make
## it the active/selected sheet - mirroring the author's commit.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new worksheet after the last existing sheet ----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "RTM"

# --- 2. Write header row -------------------------------------------------
$headers = @(
    "Requirement_ID",
    "Requirement_Description",
    "Reconciliation_Rule_ID",
    "Test_Case_ID",
    "Test_Case_Status",
    "Defect_ID",
    "Remarks"
)
for ($c = 0; $c -lt $headers.Length; $c++) {
    $newSheet.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# --- 3. Write the requirement traceability rows --------------------------
$data = @(
    @("RQ_01", "Only PROCESSED transactions from Source A should be reconciled", "RR_01", "TC_01", "PASS", "–", "FAILED transactions correctly excluded"),
    @("RQ_02", "Transaction IDs must be unique across both sources", "RR_02", "TC_02", "PASS", "–", "No duplicate transaction IDs found"),
    @("RQ_03", "Eligible transactions must exist in both Source A and Ledger", "RR_03", "TC_03", "PASS", "DEF_003, DEF_004", "Extra ledger records identified"),
    @("RQ_04", "Transaction amount must match between Source A and Ledger", "RR_04", "TC_06", "PASS", "DEF_001, DEF_002", "Amount mismatches correctly flagged"),
    @("RQ_05", "Currency should be consistent across both sources", "RR_05", "TC_08", "PASS", "–", "Currency matched for all eligible records"),
    @("RQ_06", "FAILED or CANCELLED transactions must be excluded from reconciliation", "RR_06", "TC_05", "PASS", "–", "Exclusion validated successfully"),
    @("RQ_07", "Reconciliation results must classify records correctly", "RR_07", "TC_07", "PASS", "–", "MATCHED / MISMATCH / EXTRA identified"),
    @("RQ_08", "Reconciliation report must be generated for each batch cycle", "RR_08", "TC_09", "PASS", "–", "Report generated successfully")
)
for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt 7; $c++) {
        $newSheet.Cells.Item($r + 2, $c + 1).Value = $data[$r][$c]
    }
}

# --- 4. Bold the header row (captured into the table's header style) ----
$newSheet.Range("A1:G1").Font.Bold = $true

# --- 5. Turn the range into a native Excel Table (ListObject) -----------
$tbl = $newSheet.ListObjects.Add(1, $newSheet.Range("A1:G9"), $null, 1)
$tbl.Name = "Table6"
$tbl.TableStyle = "TableStyleMedium2"

# --- 6. Column widths (matching the authored layout) ---------------------
$newSheet.Columns.Item(1).ColumnWidth = 15.966796875
$newSheet.Columns.Item(2).ColumnWidth = 60.25
$newSheet.Columns.Item(3).ColumnWidth = 21.67578125
$newSheet.Columns.Item(4).ColumnWidth = 13.9609375
$newSheet.Columns.Item(5).ColumnWidth = 17.509765625
$newSheet.Columns.Item(6).ColumnWidth = 15.658203125
$newSheet.Columns.Item(7).ColumnWidth = 34.01953125

# --- 7. Make RTM the active sheet with G39 selected ----------------------
$newSheet.Activate()
$newSheet.Range("G39").Select() | Out-Null
